# Auto-generated Excel COM-interop script to apply numeric updates
# to the Aegis_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1251.2106
$ws.Range("I15").Value = 1251.2106
$ws.Range("K15").Value = 3753.6318
$ws.Range("M15").Value = -3584.6318

$ws.Range("H86").Value = 43290.207
$ws.Range("I86").Value = 51726.6
$ws.Range("J86").Value = 1108.25
$ws.Range("K86").Value = 51726.6
$ws.Range("L86").Value = 1108.25
$ws.Range("M86").Value = -50603.6
$ws.Range("N86").Value = -3354.25

$ws.Range("H89").Value = 43290.207
$ws.Range("I89").Value = 51726.6
$ws.Range("J89").Value = 1108.25
$ws.Range("K89").Value = 258633
$ws.Range("L89").Value = 5541.25
$ws.Range("M89").Value = -253017
$ws.Range("N89").Value = -16773.25

$ws.Range("H107").Value = 352.66666
$ws.Range("I107").Value = 377.6842
$ws.Range("J107").Value = 115
$ws.Range("K107").Value = 377.6842
$ws.Range("L107").Value = 115
$ws.Range("M107").Value = 1542.3158
$ws.Range("N107").Value = -3955

$ws.Range("H125").Value = 1921.6666
$ws.Range("I125").Value = 2337.5715
$ws.Range("J125").Value = 1750.4117
$ws.Range("K125").Value = 21038.1435
$ws.Range("L125").Value = 15753.7053
$ws.Range("M125").Value = -18578.1435
$ws.Range("N125").Value = -20673.7053

$ws.Range("H132").Value = 17871570
$ws.Range("I132").Value = 17871570
$ws.Range("K132").Value = 53614710
$ws.Range("M132").Value = -53612180

$ws.Range("H135").Value = 2031.9048
$ws.Range("I135").Value = 1117.2222
$ws.Range("J135").Value = 2717.9167
$ws.Range("K135").Value = 10054.9998
$ws.Range("L135").Value = 24461.2503
$ws.Range("M135").Value = -7519.9998
$ws.Range("N135").Value = -29531.2503

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 67823.39999999999
$ws.Range("I2").Value = 1244.75
$ws.Range("K2").Value = 1244.75
$ws.Range("M2").Value = -1131.75

$ws.Range("H74").Value = 733.8889
$ws.Range("I74").Value = 554.2857
$ws.Range("J74").Value = 1362.5
$ws.Range("K74").Value = 554.2857
$ws.Range("L74").Value = 1362.5
$ws.Range("M74").Value = 319.7143
$ws.Range("N74").Value = -3110.5

$ws.Range("H77").Value = 733.8889
$ws.Range("I77").Value = 554.2857
$ws.Range("J77").Value = 1362.5
$ws.Range("K77").Value = 2771.4285
$ws.Range("L77").Value = 6812.5
$ws.Range("M77").Value = 1596.5715
$ws.Range("N77").Value = -15548.5

$ws.Range("H102").Value = 52203
$ws.Range("I102").Value = 101987.4
$ws.Range("J102").Value = 2418.6
$ws.Range("K102").Value = 101987.4
$ws.Range("L102").Value = 2418.6
$ws.Range("M102").Value = -100365.4
$ws.Range("N102").Value = -5662.6

$ws.Range("H116").Value = 67823.39999999999
$ws.Range("I116").Value = 1244.75
$ws.Range("K116").Value = 1244.75
$ws.Range("M116").Value = 1049.25

$ws.Range("H122").Value = 1429.963
$ws.Range("I122").Value = 1411.45
$ws.Range("K122").Value = 4234.35
$ws.Range("M122").Value = -1784.35

$ws.Range("H132").Value = 4597.8
$ws.Range("I132").Value = 4882.7646
$ws.Range("J132").Value = 3992.25
$ws.Range("K132").Value = 14648.2938
$ws.Range("L132").Value = 11976.75
$ws.Range("M132").Value = -12118.2938
$ws.Range("N132").Value = -17036.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 67823.39999999999
$ws.Range("I3").Value = 1244.75
$ws.Range("K3").Value = 1244.75
$ws.Range("M3").Value = -1130.75

$ws.Range("H86").Value = 48421.668
$ws.Range("I86").Value = 87939.08
$ws.Range("J86").Value = 1719.2727
$ws.Range("K86").Value = 87939.08
$ws.Range("L86").Value = 1719.2727
$ws.Range("M86").Value = -86816.08
$ws.Range("N86").Value = -3965.2727

$ws.Range("H89").Value = 48421.668
$ws.Range("I89").Value = 87939.08
$ws.Range("J89").Value = 1719.2727
$ws.Range("K89").Value = 439695.4
$ws.Range("L89").Value = 8596.363499999999
$ws.Range("M89").Value = -434079.4
$ws.Range("N89").Value = -19828.3635

$ws.Range("H94").Value = 472.64285
$ws.Range("I94").Value = 436.47058
$ws.Range("J94").Value = 528.5454999999999
$ws.Range("K94").Value = 436.47058
$ws.Range("L94").Value = 528.5454999999999
$ws.Range("M94").Value = 14.52942000000002
$ws.Range("N94").Value = -1430.5455

$ws.Range("H134").Value = 2095.1875
$ws.Range("I134").Value = 2228.1924
$ws.Range("J134").Value = 1518.8334
$ws.Range("K134").Value = 6684.5772
$ws.Range("L134").Value = 4556.5002
$ws.Range("M134").Value = -4149.5772
$ws.Range("N134").Value = -9626.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 100006270
$ws.Range("I132").Value = 142865920
$ws.Range("K132").Value = 428597760
$ws.Range("M132").Value = -428595230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 91005050
$ws.Range("I80").Value = 200207400
$ws.Range("J80").Value = 3093.3333
$ws.Range("K80").Value = 200207400
$ws.Range("L80").Value = 3093.3333
$ws.Range("M80").Value = -200206402
$ws.Range("N80").Value = -5089.3333

$ws.Range("H83").Value = 91005050
$ws.Range("I83").Value = 200207400
$ws.Range("J83").Value = 3093.3333
$ws.Range("K83").Value = 1001037000
$ws.Range("L83").Value = 15466.6665
$ws.Range("M83").Value = -1001032008
$ws.Range("N83").Value = -25450.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 40049.08
$ws.Range("I40").Value = 72394.92999999999
$ws.Range("K40").Value = 72394.92999999999
$ws.Range("M40").Value = -72258.92999999999

$ws.Range("H61").Value = 1830.3889
$ws.Range("I61").Value = 1719
$ws.Range("K61").Value = 1719
$ws.Range("M61").Value = -1517

$ws.Range("H113").Value = 1830.3889
$ws.Range("I113").Value = 1719
$ws.Range("K113").Value = 1719
$ws.Range("M113").Value = 451

$ws.Range("H132").Value = 4608.773
$ws.Range("I132").Value = 4605.278
$ws.Range("J132").Value = 4624.5
$ws.Range("K132").Value = 13815.834
$ws.Range("L132").Value = 13873.5
$ws.Range("M132").Value = -11285.834
$ws.Range("N132").Value = -18933.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 932
$ws.Range("I113").Value = 799.75
$ws.Range("K113").Value = 2399.25
$ws.Range("M113").Value = -229.25

$ws.Range("H132").Value = 2288.3
$ws.Range("I132").Value = 2309.375
$ws.Range("J132").Value = 2204
$ws.Range("K132").Value = 6928.125
$ws.Range("L132").Value = 6612
$ws.Range("M132").Value = -4398.125
